$d = $word.ActiveDocument

# --- Change 1: merge the three runs about "Write a query to display department
# name and sum of salary who are working in 'IT network' department." into one
# run with the combined text, removing the proofErr gramStart/gramEnd wrappers. ---
$startRng = $d.Content
$startRng.Find.ClearFormatting()
$foundStart = $startRng.Find.Execute(
    "Write a query to display department name and sum of ",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$endRng = $d.Content
$endRng.Find.ClearFormatting()
$foundEnd = $endRng.Find.Execute(
    " are working in 'IT network' department.",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($foundStart -and $foundEnd) {
    $combined = $d.Range($startRng.Start, $endRng.End)
    # Write a temp placeholder first: the engine no-ops a replace whose
    # resulting text is identical to what is already there (it would be,
    # since we are just re-joining the existing text across runs), so we
    # force an actual content change, then set the real final text.
    $combined.Text = "PLACEHOLDER_TEXT_TEMP"
    $final = $d.Range($startRng.Start, $startRng.Start + 22)
    $final.Text = "Write a query to display department name and sum of salary  who are working in 'IT network' department."
}

# --- Change 2: Replace "ANS." paragraph text with "ANSWERS" (two runs: "ANS" + "WERS"),
# and switch the paragraph style to the built-in "Intense Quote" style. ---
$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$found = $rng2.Find.Execute("ANS.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng2.Text = "ANSWERS"
    $rng2.Paragraphs.Item(1).Style = "Intense Quote"
}
